$d = $word.ActiveDocument

$d.Content.Find.Execute("38×18=", $true, $false, $false, $false, $false, $true, 1, $false, "42×99=", 2)
$d.Content.Find.Execute("71×15=", $true, $false, $false, $false, $false, $true, 1, $false, "36×56=", 2)
$d.Content.Find.Execute("24×68=", $true, $false, $false, $false, $false, $true, 1, $false, "97×16=", 2)
$d.Content.Find.Execute("98×65=", $true, $false, $false, $false, $false, $true, 1, $false, "16×93=", 2)
$d.Content.Find.Execute("22×64=", $true, $false, $false, $false, $false, $true, 1, $false, "89×97=", 2)
$d.Content.Find.Execute("50×82=", $true, $false, $false, $false, $false, $true, 1, $false, "92×11=", 2)
$d.Content.Find.Execute("17×35=", $true, $false, $false, $false, $false, $true, 1, $false, "41×12=", 2)
$d.Content.Find.Execute("71×94=", $true, $false, $false, $false, $false, $true, 1, $false, "38×73=", 2)
$d.Content.Find.Execute("24×95=", $true, $false, $false, $false, $false, $true, 1, $false, "38×41=", 2)
$d.Content.Find.Execute("53×22=", $true, $false, $false, $false, $false, $true, 1, $false, "28×65=", 2)
$d.Content.Find.Execute("21×43=", $true, $false, $false, $false, $false, $true, 1, $false, "16×72=", 2)
$d.Content.Find.Execute("77×91=", $true, $false, $false, $false, $false, $true, 1, $false, "12×49=", 2)
$d.Content.Find.Execute("81×25=", $true, $false, $false, $false, $false, $true, 1, $false, "93×29=", 2)
$d.Content.Find.Execute("35×56=", $true, $false, $false, $false, $false, $true, 1, $false, "81×95=", 2)
$d.Content.Find.Execute("57×12=", $true, $false, $false, $false, $false, $true, 1, $false, "13×57=", 2)
$d.Content.Find.Execute("11×84=", $true, $false, $false, $false, $false, $true, 1, $false, "33×84=", 2)
$d.Content.Find.Execute("30×78=", $true, $false, $false, $false, $false, $true, 1, $false, "93×27=", 2)
$d.Content.Find.Execute("53×87=", $true, $false, $false, $false, $false, $true, 1, $false, "56×67=", 2)
$d.Content.Find.Execute("79×97=", $true, $false, $false, $false, $false, $true, 1, $false, "70×24=", 2)
$d.Content.Find.Execute("58×87=", $true, $false, $false, $false, $false, $true, 1, $false, "92×50=", 2)
$d.Content.Find.Execute("32×29=", $true, $false, $false, $false, $false, $true, 1, $false, "96×84=", 2)
$d.Content.Find.Execute("48×83=", $true, $false, $false, $false, $false, $true, 1, $false, "50×38=", 2)
$d.Content.Find.Execute("82×56=", $true, $false, $false, $false, $false, $true, 1, $false, "60×71=", 2)
$d.Content.Find.Execute("73×78=", $true, $false, $false, $false, $false, $true, 1, $false, "43×39=", 2)
$d.Content.Find.Execute("58×91=", $true, $false, $false, $false, $false, $true, 1, $false, "18×16=", 2)
